# Applies the cryptos list update described by the commit diff.
# Each cell is set to its new literal text value. Numeric-looking values
# in column D are prefixed with a leading apostrophe so Excel keeps them
# as text (matching the workbook's original inlineStr/text storage)
# instead of silently reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "92.539.32"
$ws.Range("E2").Value = "  +2.50%  "

# Row 3
$ws.Range("D3").Value = "3.134.39"
$ws.Range("E3").Value = "  +1.75%  "

# Row 4
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'245.27"
$ws.Range("E5").Value = "  +1.20%  "

# Row 6
$ws.Range("D6").Value = "'620.83"
$ws.Range("E6").Value = "  +0.80%  "

# Row 7
$ws.Range("E7").Value = "  -1.55%  "

# Row 8
$ws.Range("D8").Value = "'0.394"
$ws.Range("E8").Value = "  +8.63%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10
$ws.Range("D10").Value = "3.134.12"
$ws.Range("E10").Value = "  +1.82%  "

# Row 11
$ws.Range("D11").Value = "'0.748"
$ws.Range("E11").Value = "  +2.68%  "

# Row 12
$ws.Range("D12").Value = "'0.204"
$ws.Range("E12").Value = "  +1.23%  "

# Row 13
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  +4.23%  "

# Row 14
$ws.Range("D14").Value = "'35.18"
$ws.Range("E14").Value = "  +1.53%  "

# Row 15
$ws.Range("D15").Value = "'5.57"
$ws.Range("E15").Value = "  +2.17%  "

# Row 16
$ws.Range("D16").Value = "92.187.00"
$ws.Range("E16").Value = "  +2.06%  "

# Row 17
$ws.Range("D17").Value = "3.718.40"
$ws.Range("E17").Value = "  +1.63%  "

# Row 18
$ws.Range("D18").Value = "3.149.25"
$ws.Range("E18").Value = "  +2.11%  "

# Row 19
$ws.Range("D19").Value = "'3.73"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20
$ws.Range("D20").Value = "'15.01"
$ws.Range("E20").Value = "  +3.88%  "

# Row 21
$ws.Range("D21").Value = "'5.89"
$ws.Range("E21").Value = "  +3.02%  "

# Row 22
$ws.Range("D22").Value = "'9.50"
$ws.Range("E22").Value = "  +6.04%  "

# Row 23
$ws.Range("D23").Value = "'453.13"
$ws.Range("E23").Value = "  +3.42%  "

# Row 24
$ws.Range("D24").Value = "'0.0000204"
$ws.Range("E24").Value = "  -1.95%  "

# Row 25
$ws.Range("D25").Value = "'5.93"
$ws.Range("E25").Value = "  +6.16%  "

# Row 26
$ws.Range("D26").Value = "'88.57"
$ws.Range("E26").Value = "  -2.58%  "

# Row 27
$ws.Range("D27").Value = "'11.85"
$ws.Range("E27").Value = "  +0.60%  "

# Row 28
$ws.Range("D28").Value = "3.304.81"
$ws.Range("E28").Value = "  +1.90%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  +26.95%  "

# Row 31
$ws.Range("D31").Value = "'0.235"
$ws.Range("E31").Value = "  -1.28%  "

# Row 32
$ws.Range("D32").Value = "'0.169"
$ws.Range("E32").Value = "  -5.13%  "

# Row 33
$ws.Range("D33").Value = "'9.45"
$ws.Range("E33").Value = "  +3.65%  "

# Row 34
$ws.Range("D34").Value = "'0.177"
$ws.Range("E34").Value = "  +5.94%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +4.72%  "

# Row 36
$ws.Range("D36").Value = "'8.16"
$ws.Range("E36").Value = "  +7.93%  "

# Row 37
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'26.62"
$ws.Range("E37").Value = "  +1.53%  "

# Row 38
$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").Value = "'4.25"
$ws.Range("E38").Value = "  +1.19%  "

# Row 39
$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +3.30%  "

# Row 40
$ws.Range("D40").Value = "'495.76"
$ws.Range("E40").Value = "  +2.03%  "

# Row 41
$ws.Range("E41").Value = "  +3.90%  "

# Row 42
$ws.Range("D42").Value = "'0.443"
$ws.Range("E42").Value = "  +6.47%  "

# Row 43
$ws.Range("D43").Value = "'3.50"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44
$ws.Range("D44").Value = "'22.23"
$ws.Range("E44").Value = "  +0.27%  "

# Row 46
$ws.Range("D46").Value = "'1.96"
$ws.Range("E46").Value = "  +4.31%  "

# Row 47
$ws.Range("D47").Value = "'159.18"
$ws.Range("E47").Value = "  +3.56%  "

# Row 48
$ws.Range("D48").Value = "'0.707"
$ws.Range("E48").Value = "  +3.93%  "

# Row 49
$ws.Range("E49").Value = "  +3.84%  "

# Row 50
$ws.Range("D50").Value = "'0.0333"
$ws.Range("E50").Value = "  +8.21%  "

# Row 51
$ws.Range("D51").Value = "'4.43"
$ws.Range("E51").Value = "  +0.64%  "
